$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H116").Value = 43499.055
$ws.Range("J116").Value = 44999
$ws.Range("L116").Value = 44999
$ws.Range("N116").Value = -51883
$ws.Range("H135").Value = 706.86206
$ws.Range("I135").Value = 412.52173
$ws.Range("K135").Value = 3712.69557
$ws.Range("M135").Value = -1177.69557
$ws.Range("H137").Value = 4116.881
$ws.Range("I137").Value = 5603.4585
$ws.Range("J137").Value = 2134.7778
$ws.Range("K137").Value = 16810.3755
$ws.Range("L137").Value = 6404.3334
$ws.Range("M137").Value = -14260.3755
$ws.Range("N137").Value = -11504.3334

$ws = $wb.Worksheets("ARM")
$ws.Range("H32").Value = 2855.3923
$ws.Range("I32").Value = 3119.4775
$ws.Range("J32").Value = 1380.9166
$ws.Range("K32").Value = 3119.4775
$ws.Range("L32").Value = 1380.9166
$ws.Range("M32").Value = -2832.4775
$ws.Range("N32").Value = -1954.9166
$ws.Range("H45").Value = 7012.9287
$ws.Range("I45").Value = 8820.611000000001
$ws.Range("K45").Value = 8820.611000000001
$ws.Range("M45").Value = -8443.611000000001
$ws.Range("H74").Value = 5229.644
$ws.Range("J74").Value = 6983.8335
$ws.Range("L74").Value = 6983.8335
$ws.Range("N74").Value = -8731.833500000001
$ws.Range("H77").Value = 5229.644
$ws.Range("J77").Value = 6983.8335
$ws.Range("L77").Value = 34919.1675
$ws.Range("N77").Value = -43655.1675
$ws.Range("H135").Value = 85999
$ws.Range("J135").Value = 85999
$ws.Range("L135").Value = 85999
$ws.Range("N135").Value = -96139

$ws = $wb.Worksheets("BSM")
$ws.Range("H22").Value = 296.92856
$ws.Range("I22").Value = 262.07693
$ws.Range("J22").Value = 750
$ws.Range("K22").Value = 262.07693
$ws.Range("L22").Value = 750
$ws.Range("M22").Value = -89.07693
$ws.Range("N22").Value = -1096

$ws = $wb.Worksheets("CRP")
$ws.Range("H22").Value = 1953.2174
$ws.Range("I22").Value = 1924
$ws.Range("K22").Value = 1924
$ws.Range("M22").Value = -1574
$ws.Range("H31").Value = 2162.2778
$ws.Range("I31").Value = 1658.2727
$ws.Range("J31").Value = 2954.2856
$ws.Range("K31").Value = 1658.2727
$ws.Range("L31").Value = 2954.2856
$ws.Range("M31").Value = -1363.2727
$ws.Range("N31").Value = -3544.2856
$ws.Range("H34").Value = 2162.2778
$ws.Range("I34").Value = 1658.2727
$ws.Range("J34").Value = 2954.2856
$ws.Range("K34").Value = 1658.2727
$ws.Range("L34").Value = 2954.2856
$ws.Range("M34").Value = -1456.2727
$ws.Range("N34").Value = -3358.2856
$ws.Range("H58").Value = 21068.875
$ws.Range("J58").Value = 24248.5
$ws.Range("L58").Value = 24248.5
$ws.Range("N58").Value = -24654.5
$ws.Range("H99").Value = 9569.388999999999
$ws.Range("I99").Value = 7139.846
$ws.Range("J99").Value = 10942.608
$ws.Range("K99").Value = 7139.846
$ws.Range("L99").Value = 10942.608
$ws.Range("M99").Value = -5641.846
$ws.Range("N99").Value = -13938.608
$ws.Range("H107").Value = 1488.1666
$ws.Range("I107").Value = 1155.3
$ws.Range("K107").Value = 1155.3
$ws.Range("M107").Value = 764.7
$ws.Range("H126").Value = 9569.388999999999
$ws.Range("I126").Value = 7139.846
$ws.Range("J126").Value = 10942.608
$ws.Range("K126").Value = 21419.538
$ws.Range("L126").Value = 32827.824
$ws.Range("M126").Value = -18949.538
$ws.Range("N126").Value = -37767.824
$ws.Range("H132").Value = 7834.0625
$ws.Range("I132").Value = 3021.4285
$ws.Range("K132").Value = 9064.2855
$ws.Range("M132").Value = -6534.2855
$ws.Range("H134").Value = 2208.2334
$ws.Range("I134").Value = 1925
$ws.Range("J134").Value = 2697.4546
$ws.Range("K134").Value = 5775
$ws.Range("L134").Value = 8092.3638
$ws.Range("M134").Value = -3240
$ws.Range("N134").Value = -13162.3638
$ws.Range("H136").Value = 21068.875
$ws.Range("J136").Value = 24248.5
$ws.Range("L136").Value = 72745.5
$ws.Range("N136").Value = -77845.5

$ws = $wb.Worksheets("CUL")
$ws.Range("H5").Value = 839.125
$ws.Range("I5").Value = 602.9286
$ws.Range("J5").Value = 1169.8
$ws.Range("K5").Value = 1808.7858
$ws.Range("L5").Value = 3509.4
$ws.Range("M5").Value = -1696.7858
$ws.Range("N5").Value = -3733.4
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("H97").Value = 467.55554
$ws.Range("J97").Value = 474.2857
$ws.Range("L97").Value = 1422.8571
$ws.Range("N97").Value = -2414.8571
$ws.Range("H107").Value = 507.46155
$ws.Range("I107").Value = 532.25
$ws.Range("J107").Value = 496.44446
$ws.Range("K107").Value = 1596.75
$ws.Range("L107").Value = 1489.33338
$ws.Range("M107").Value = 323.25
$ws.Range("N107").Value = -5329.33338
$ws.Range("H129").Value = 1842.9445
$ws.Range("I129").Value = 729.63635
$ws.Range("K129").Value = 2188.90905
$ws.Range("M129").Value = 2811.09095
$ws.Range("H135").Value = 839.125
$ws.Range("I135").Value = 602.9286
$ws.Range("J135").Value = 1169.8
$ws.Range("K135").Value = 5426.3574
$ws.Range("L135").Value = 10528.2
$ws.Range("M135").Value = -2891.3574
$ws.Range("N135").Value = -15598.2
$ws.Range("M87").ClearContents()
$ws.Range("M90").ClearContents()

$ws = $wb.Worksheets("GSM")
$ws.Range("H70").Value = 5882.1875
$ws.Range("I70").Value = 5688.625
$ws.Range("J70").Value = 6075.75
$ws.Range("K70").Value = 5688.625
$ws.Range("L70").Value = 6075.75
$ws.Range("M70").Value = -5418.625
$ws.Range("N70").Value = -6615.75
$ws.Range("H73").Value = 5882.1875
$ws.Range("I73").Value = 5688.625
$ws.Range("J73").Value = 6075.75
$ws.Range("K73").Value = 5688.625
$ws.Range("L73").Value = 6075.75
$ws.Range("M73").Value = -4752.625
$ws.Range("N73").Value = -7947.75
$ws.Range("H80").Value = 28402760
$ws.Range("I80").Value = 50716964
$ws.Range("J80").Value = 2867.2727
$ws.Range("K80").Value = 50716964
$ws.Range("L80").Value = 2867.2727
$ws.Range("M80").Value = -50715966
$ws.Range("N80").Value = -4863.2727
$ws.Range("H83").Value = 28402760
$ws.Range("I83").Value = 50716964
$ws.Range("J83").Value = 2867.2727
$ws.Range("K83").Value = 253584820
$ws.Range("L83").Value = 14336.3635
$ws.Range("M83").Value = -253579828
$ws.Range("N83").Value = -24320.3635
$ws.Range("H122").Value = 2672.7368
$ws.Range("I122").Value = 1963
$ws.Range("K122").Value = 5889
$ws.Range("M122").Value = -3439
$ws.Range("H132").Value = 6110.706
$ws.Range("I132").Value = 5532.25
$ws.Range("K132").Value = 16596.75
$ws.Range("M132").Value = -14066.75

$ws = $wb.Worksheets("LTW")
$ws.Range("H93").Value = 5883431
$ws.Range("I93").Value = 6452489.5
$ws.Range("J93").Value = 3160
$ws.Range("K93").Value = 6452489.5
$ws.Range("L93").Value = 3160
$ws.Range("M93").Value = -6451241.5
$ws.Range("N93").Value = -5656
$ws.Range("H132").Value = 32218
$ws.Range("I132").Value = 37450.695
$ws.Range("K132").Value = 112352.085
$ws.Range("M132").Value = -109822.085
$ws.Range("H136").Value = 23852.787
$ws.Range("I136").Value = 24157.26
$ws.Range("J136").Value = 22482.666
$ws.Range("K136").Value = 72471.78
$ws.Range("L136").Value = 67447.998
$ws.Range("M136").Value = -69921.78
$ws.Range("N136").Value = -72547.99800000001

$ws = $wb.Worksheets("WVR")
$ws.Range("H107").Value = 1799.5714
$ws.Range("J107").Value = 500
$ws.Range("L107").Value = 1500
$ws.Range("N107").Value = -5340
$ws.Range("H126").Value = 5423.3335
$ws.Range("I126").Value = 5081.2
$ws.Range("J126").Value = 7134
$ws.Range("K126").Value = 15243.6
$ws.Range("L126").Value = 21402
$ws.Range("M126").Value = -12773.6
$ws.Range("N126").Value = -26342
